$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (K2:T2)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4040993333333333
$ws.Range("N2").Value = 1.212298
$ws.Range("O2").Value = 0.1080113302049822
$ws.Range("P2").Value = 0.1080113302049822
$ws.Range("Q2").Value = 0.02723036827644445
$ws.Range("R2").Value = 0.245073314488
$ws.Range("S2").Value = 0.1080113302049822
$ws.Range("T2").Value = 0.1080113302049822

# Update row 3 values (M3:T3)
$ws.Range("M3").Value = 3.337168666666667
$ws.Range("N3").Value = 10.011506
$ws.Range("O3").Value = 0.8919886697950178
$ws.Range("P3").Value = 0.8919886697950178
$ws.Range("Q3").Value = 0.2248762229928889
$ws.Range("R3").Value = 2.023886006936
$ws.Range("S3").Value = 0.8919886697950178
$ws.Range("T3").Value = 0.8919886697950178

# Delete row 4 entirely (Resolving-Mac row)
$ws.Range("A4:T4").EntireRow.Delete()
